try {
    $wb = $excel.ActiveWorkbook

    # ------------------------------------------------------------------
    # 1. Features sheet: rename "chemical_formula" header to "chem_formula"
    #    (molecular_weight column keeps its text, just shifts string id)
    # ------------------------------------------------------------------
    $wsFeatures = $wb.Worksheets.Item("Features")
    $wsFeatures.Range("E1").Value2 = "chem_formula"

    # Update cursor/selection to E1 on the Features sheet (matches diff)
    $wsFeatures.Activate()
    $wsFeatures.Range("E1").Select()

    # ------------------------------------------------------------------
    # 2. ISTDs sheet: insert a new "istd_conc_ngml" column (C1), pushing
    #    the existing "remarks" column (was C1) to D1. Give the new
    #    header the same "required" (blue) formatting as the other
    #    required headers (A1 / B1).
    # ------------------------------------------------------------------
    $wsIstds = $wb.Worksheets.Item("ISTDs")

    # Move remarks header from C1 to D1 first
    $wsIstds.Range("D1").Value2 = $wsIstds.Range("C1").Value2

    # Write the new header into C1
    $wsIstds.Range("C1").Value2 = "istd_conc_ngml"

    # Copy the "required" header formatting (blue font) from B1 onto C1
    $wsIstds.Range("B1").Copy()
    $wsIstds.Range("C1").PasteSpecial(-4122)  # xlPasteFormats
    $excel.CutCopyMode = 0

    # Add comments describing the unit columns
    $wsIstds.Range("B1").AddComment("Define either nmolar or ngmL")
    $wsIstds.Range("C1").AddComment("Define either ngmL or nmolar")

    # Update cursor/selection to B1 on the ISTDs sheet (matches diff)
    $wsIstds.Activate()
    $wsIstds.Range("B1").Select()

    # ------------------------------------------------------------------
    # 3. QCconcentrations sheet: selection/cursor moved to D13 (no data
    #    change on this sheet).
    # ------------------------------------------------------------------
    $wsQc = $wb.Worksheets.Item("QCconcentrations")
    $wsQc.Activate()
    $wsQc.Range("D13").Select()

    # ------------------------------------------------------------------
    # Restore Features as the active/selected sheet (it was tabSelected
    # in the original workbook).
    # ------------------------------------------------------------------
    $wsFeatures.Activate()

    Write-Output "done"
} catch {
    Write-Output "ERROR: $_"
}
